# feat: add 2022-Q3 data
#
# 1) Insert a new worksheet "2022-Q3" right before "2022-Q2" and populate it
#    with the quarterly fund-holding data.
# 2) Insert a new row at the top of the "总计" (summary) sheet's data table
#    for the 2022-Q3 aggregate figures, shifting the older rows down.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Force a numeric-looking string to be stored as text (keeps leading
    # zeros / decimal formatting exactly as supplied) by using Excel's
    # standard "leading apostrophe" text-prefix convention.
    $range.Value = "'" + $text
}

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet
# ---------------------------------------------------------------------------
$ws2022Q2 = $wb.Worksheets.Item("2022-Q2")
$wsQ3 = $wb.Worksheets.Add($ws2022Q2)
$wsQ3.Name = "2022-Q3"

# Header row (bold, bordered, centered horizontally, top-aligned vertically
# - matches the look of the other quarter sheets)
$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

$headerRange = $wsQ3.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows: index, fund code, fund name, size, stock position, position
# ratio, held market value (100M yuan), position rank
$dataQ3 = @(
    @(0, "159993", "鹏华国证证券龙头ETF",                               "13.39", "98.06", "4.07", "0.5450", 9),
    @(1, "016935", "景顺长城中证500指数增强C",                          "15.57", "93.89", "2.33", "0.3628", 1),
    @(2, "000978", "景顺长城量化精选股票",                               "7.14",  "93.64", "2.18", "0.1557", 1),
    @(3, "515760", "华夏中证浙江国资创新发展ETF",                        "2.04",  "99.57", "6.38", "0.1302", 3),
    @(4, "008851", "景顺长城量化对冲策略三个月定期开放灵活配置混合",         "2.96",  "64.77", "1.51", "0.0447", 1),
    @(5, "015859", "宝盈国证证券龙头指数A",                              "0.12",  "94.13", "3.92", "0.0047", 9),
    @(6, "006611", "人保中证500指数",                                   "0.39",  "92.75", "0.72", "0.0028", 3),
    @(7, "015860", "宝盈国证证券龙头指数C",                              "0.05",  "94.13", "3.92", "0.0020", 9),
    @(8, "006682", "景顺长城中证500指数增强A",                          "0.00",  "93.89", "2.33", "__NUM0__", 1)
)

$rowIdx = 2
foreach ($row in $dataQ3) {
    $aCell = $wsQ3.Range("A$rowIdx")
    $aCell.Value = $row[0]
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    Set-TextValue $wsQ3.Range("B$rowIdx") $row[1]
    $wsQ3.Range("C$rowIdx").Value = $row[2]
    Set-TextValue $wsQ3.Range("D$rowIdx") $row[3]
    Set-TextValue $wsQ3.Range("E$rowIdx") $row[4]
    Set-TextValue $wsQ3.Range("F$rowIdx") $row[5]

    if ($row[6] -eq "__NUM0__") {
        $wsQ3.Range("G$rowIdx").Value = 0
    } else {
        Set-TextValue $wsQ3.Range("G$rowIdx") $row[6]
    }

    $wsQ3.Range("H$rowIdx").Value = $row[7]
    $rowIdx++
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row for 2022-Q3
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()
# The inserted row inherits formatting from the row above it; clear it so
# the new data cells start out unstyled, like the existing data rows.
$wsTotal.Range("B2:D2").ClearFormats()

$a2 = $wsTotal.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 9
$wsTotal.Range("D2").Value = 1.25

# Renumber the A column (sequential index counter) for the rows that were
# pushed down by the insertion.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("A6").Value = 4
$wsTotal.Range("A7").Value = 5
$wsTotal.Range("A8").Value = 6
$wsTotal.Range("A9").Value = 7

# Restore the "总计" sheet as the active sheet (matches original workbook
# state; the newly-added sheet should not remain the selected tab).
$wsTotal.Activate()
$wsTotal.Range("A1").Select() | Out-Null

Write-Host "2022-Q3 data added successfully"
